$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8457584381103516
$ws.Range("B1").Value = 1.311477065086365
$ws.Range("C1").Value = 4.653108596801758
$ws.Range("D1").Value = 3.895499467849731
$ws.Range("E1").Value = 0.5210281014442444
